$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 7786
$ws.Range("F13").Value = 5730
$ws.Range("F15").Value = 2689
$ws.Range("F16").Value = 1189
$ws.Range("F21").Value = 584
$ws.Range("F22").Value = 13
$ws.Range("F23").Value = 3757
$ws.Range("F25").Value = 51
$ws.Range("F29").Value = 3912
$ws.Range("F31").Value = 60
$ws.Range("F33").Value = 376
$ws.Range("F35").Value = 362
$ws.Range("F36").Value = 1359
$ws.Range("F40").Value = 3159
$ws.Range("F41").Value = 61
$ws.Range("F44").Value = 3334
$ws.Range("F46").Value = 2315
$ws.Range("F47").Value = 14

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 62
$ws.Range("F10").Value = 5

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1354

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1354
$ws.Range("F5").Value = 7786
$ws.Range("F13").Value = 5730
$ws.Range("F15").Value = 2689
$ws.Range("F16").Value = 1189
$ws.Range("F23").Value = 584
$ws.Range("F25").Value = 3757
$ws.Range("F27").Value = 51
$ws.Range("F30").Value = 3913
$ws.Range("F31").Value = 60
$ws.Range("F32").Value = 376
$ws.Range("F34").Value = 362
$ws.Range("F35").Value = 62
$ws.Range("F36").Value = 1359
$ws.Range("F41").Value = 3159
$ws.Range("F42").Value = 61
$ws.Range("F45").Value = 3334
$ws.Range("F47").Value = 2315
